$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix H1257 / H1258: IsReleaseDay True -> (no change needed except these two releases became release days) ---
$ws.Range("H1257").Value = "'True"
$ws.Range("H1258").Value = "'True"

# Row 1259
$ws.Range("A1259").Value = '2024-11-03'
$ws.Range("B1259").Value = '17:30'
$ws.Range("C1259").Value = '18:15'
$ws.Range("D1259").Value = '0h 45m'
$ws.Range("E1259").Value = '#maintenance'
$ws.Range("G1259").Value = "'False"
$ws.Range("H1259").Value = "'False"

# Row 1260
$ws.Range("A1260").Value = '2024-11-03'
$ws.Range("B1260").Value = '21:00'
$ws.Range("C1260").Value = '00:00'
$ws.Range("D1260").Value = '3h 00m'
$ws.Range("E1260").Value = '#python'
$ws.Range("F1260").Value = 'nwreadinglist v4.0.0'
$ws.Range("G1260").Value = "'True"
$ws.Range("H1260").Value = "'False"

# Row 1261
$ws.Range("A1261").Value = '2024-11-04'
$ws.Range("B1261").Value = '19:30'
$ws.Range("C1261").Value = '21:00'
$ws.Range("D1261").Value = '1h 30m'
$ws.Range("E1261").Value = '#maintenance'
$ws.Range("G1261").Value = "'False"
$ws.Range("H1261").Value = "'False"

# Row 1262
$ws.Range("A1262").Value = '2024-11-05'
$ws.Range("B1262").Value = '09:30'
$ws.Range("C1262").Value = '13:30'
$ws.Range("D1262").Value = '4h 00m'
$ws.Range("E1262").Value = '#studying'
$ws.Range("G1262").Value = "'False"
$ws.Range("H1262").Value = "'False"

# Row 1263
$ws.Range("A1263").Value = '2024-11-05'
$ws.Range("B1263").Value = '16:00'
$ws.Range("C1263").Value = '19:00'
$ws.Range("D1263").Value = '3h 00m'
$ws.Range("E1263").Value = '#python'
$ws.Range("F1263").Value = 'nwshared v1.7.0'
$ws.Range("G1263").Value = "'True"
$ws.Range("H1263").Value = "'True"

# Row 1264
$ws.Range("A1264").Value = '2024-11-07'
$ws.Range("B1264").Value = '08:00'
$ws.Range("C1264").Value = '08:45'
$ws.Range("D1264").Value = '0h 45m'
$ws.Range("E1264").Value = '#python'
$ws.Range("F1264").Value = 'nwreadinglist v4.0.0'
$ws.Range("G1264").Value = "'True"
$ws.Range("H1264").Value = "'False"

# Row 1265
$ws.Range("A1265").Value = '2024-11-07'
$ws.Range("B1265").Value = '17:15'
$ws.Range("C1265").Value = '17:45'
$ws.Range("D1265").Value = '0h 30m'
$ws.Range("E1265").Value = '#python'
$ws.Range("F1265").Value = 'nwreadinglist v4.0.0'
$ws.Range("G1265").Value = "'True"
$ws.Range("H1265").Value = "'False"

# Row 1266
$ws.Range("A1266").Value = '2024-11-08'
$ws.Range("B1266").Value = '08:00'
$ws.Range("C1266").Value = '08:45'
$ws.Range("D1266").Value = '0h 45m'
$ws.Range("E1266").Value = '#python'
$ws.Range("F1266").Value = 'nwreadinglist v4.0.0'
$ws.Range("G1266").Value = "'True"
$ws.Range("H1266").Value = "'False"

# Row 1267
$ws.Range("A1267").Value = '2024-11-08'
$ws.Range("B1267").Value = '17:00'
$ws.Range("C1267").Value = '17:30'
$ws.Range("D1267").Value = '0h 30m'
$ws.Range("E1267").Value = '#python'
$ws.Range("F1267").Value = 'nwreadinglist v4.0.0'
$ws.Range("G1267").Value = "'True"
$ws.Range("H1267").Value = "'False"

# Row 1268
$ws.Range("A1268").Value = '2024-11-10'
$ws.Range("B1268").Value = '10:15'
$ws.Range("C1268").Value = '15:00'
$ws.Range("D1268").Value = '4h 45m'
$ws.Range("E1268").Value = '#python'
$ws.Range("F1268").Value = 'nwreadinglist v4.0.0'
$ws.Range("G1268").Value = "'True"
$ws.Range("H1268").Value = "'False"

# Row 1269
$ws.Range("A1269").Value = '2024-11-10'
$ws.Range("B1269").Value = '15:30'
$ws.Range("C1269").Value = '19:30'
$ws.Range("D1269").Value = '4h 00m'
$ws.Range("E1269").Value = '#python'
$ws.Range("F1269").Value = 'nwreadinglist v4.0.0'
$ws.Range("G1269").Value = "'True"
$ws.Range("H1269").Value = "'False"

# Row 1270
$ws.Range("A1270").Value = '2024-11-10'
$ws.Range("B1270").Value = '20:00'
$ws.Range("C1270").Value = '21:00'
$ws.Range("D1270").Value = '1h 00m'
$ws.Range("E1270").Value = '#python'
$ws.Range("F1270").Value = 'nwreadinglist v4.0.0'
$ws.Range("G1270").Value = "'True"
$ws.Range("H1270").Value = "'False"

# Row 1271
$ws.Range("A1271").Value = '2024-11-10'
$ws.Range("B1271").Value = '22:30'
$ws.Range("C1271").Value = '01:00'
$ws.Range("D1271").Value = '2h 30m'
$ws.Range("E1271").Value = '#python'
$ws.Range("F1271").Value = 'nwreadinglist v4.0.0'
$ws.Range("G1271").Value = "'True"
$ws.Range("H1271").Value = "'False"

# Row 1272
$ws.Range("A1272").Value = '2024-11-11'
$ws.Range("B1272").Value = '10:45'
$ws.Range("C1272").Value = '13:30'
$ws.Range("D1272").Value = '2h 45m'
$ws.Range("E1272").Value = '#python'
$ws.Range("F1272").Value = 'nwreadinglist v4.0.0'
$ws.Range("G1272").Value = "'True"
$ws.Range("H1272").Value = "'True"

# Row 1273
$ws.Range("A1273").Value = '2024-11-11'
$ws.Range("B1273").Value = '14:30'
$ws.Range("C1273").Value = '20:30'
$ws.Range("D1273").Value = '6h 00m'
$ws.Range("E1273").Value = '#python'
$ws.Range("F1273").Value = 'nwreadinglist v4.0.0'
$ws.Range("G1273").Value = "'True"
$ws.Range("H1273").Value = "'True"

# Row 1274
$ws.Range("A1274").Value = '2024-11-11'
$ws.Range("B1274").Value = '20:45'
$ws.Range("C1274").Value = '22:30'
$ws.Range("D1274").Value = '1h 45m'
$ws.Range("E1274").Value = '#python'
$ws.Range("F1274").Value = 'nwreadinglist v4.0.0'
$ws.Range("G1274").Value = "'True"
$ws.Range("H1274").Value = "'True"

# Year/Month formulas (grouped to mirror original shared-formula layout)
$ws.Range("I1259").Formula = "=YEAR(A1259)"
$ws.Range("J1259").Formula = "=MONTH(A1259)"
$ws.Range("I1260:I1261").Formula = "=YEAR(A1260)"
$ws.Range("J1260:J1261").Formula = "=MONTH(A1260)"
$ws.Range("I1262:I1265").Formula = "=YEAR(A1262)"
$ws.Range("J1262:J1265").Formula = "=MONTH(A1262)"
$ws.Range("I1266").Formula = "=YEAR(A1266)"
$ws.Range("J1266").Formula = "=MONTH(A1266)"
$ws.Range("I1267:I1270").Formula = "=YEAR(A1267)"
$ws.Range("J1267:J1270").Formula = "=MONTH(A1267)"
$ws.Range("I1271:I1273").Formula = "=YEAR(A1271)"
$ws.Range("J1271:J1273").Formula = "=MONTH(A1271)"
$ws.Range("I1274").Formula = "=YEAR(A1274)"
$ws.Range("J1274").Formula = "=MONTH(A1274)"

# --- Append 16 new blank rows (1276:1291) matching the existing blank-row template (row 1275) ---
$ws.Range("A1275:J1275").Copy()
$ws.Range("A1276:J1291").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- View state: selection + frozen-pane scroll position ---
$ws.Range("H1278").Select()
$excel.ActiveWindow.ScrollRow = 1253
